$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ row=2; D='29.110.35'; E='  +0.13%  ' },
    @{ row=3; D='1.838.83'; E='  +0.05%  ' },
    @{ row=4; D='0.9981'; E='  -0.26%  ' },
    @{ row=5; D='243.38'; E='  -0.35%  ' },
    @{ row=6; D='0.6231'; E='  -1.67%  ' },
    @{ row=7; D='0.9998'; E='  -0.13%  ' },
    @{ row=8; D='0.07516'; E='  -0.98%  ' },
    @{ row=9; D='0.2944'; E='  -0.24%  ' },
    @{ row=10; D='23.38'; E='  +2.38%  ' },
    @{ row=11; D='0.07702'; E='  -0.63%  ' },
    @{ row=12; D='1.837.43'; E='  +0.06%  ' },
    @{ row=13; D='5.017'; E='  +0.26%  ' },
    @{ row=14; D='0.6763'; E='  +0.77%  ' },
    @{ row=15; D='82.91'; E='  -0.50%  ' },
    @{ row=16; D='0.000009372'; E='  -4.83%  ' },
    @{ row=17; D='5.969'; E='  -2.60%  ' },
    @{ row=18; D='29.099.04'; E='  +0.02%  ' },
    @{ row=19; D='2.075.34'; E='  -0.19%  ' },
    @{ row=20; D='12.68'; E='  +0.95%  ' },
    @{ row=21; D='222.82'; E='  -1.98%  ' },
    @{ row=22; D='1.000'; E='  -0.02%  ' },
    @{ row=23; D='7.161'; E='  -1.25%  ' },
    @{ row=24; D='1.000'; E='  +0.06%  ' },
    @{ row=25; D='160.25'; E='  -0.17%  ' },
    @{ row=26; D='0.1401'; E='  -0.65%  ' },
    @{ row=27; D='8.543'; E='  -0.13%  ' },
    @{ row=28; D='17.90'; E='  -0.51%  ' },
    @{ row=29; D='1.496'; E='  -0.45%  ' },
    @{ row=30; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.186'; E='  +1.40%  ' },
    @{ row=31; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05621'; E='  +4.40%  ' },
    @{ row=32; D='4.144'; E='  +2.38%  ' },
    @{ row=33; D='1.208'; E='  +0.36%  ' },
    @{ row=34; D='0.7493'; E='  -0.32%  ' },
    @{ row=35; D='1.849'; E='  -0.82%  ' },
    @{ row=36; D='1.147'; E='  +0.24%  ' },
    @{ row=37; D='2.664'; E='  -0.21%  ' },
    @{ row=38; D='1.236.95'; E='  -1.12%  ' },
    @{ row=39; D='2.771'; E='  +0.33%  ' },
    @{ row=40; D='0.01779'; E='  -1.11%  ' },
    @{ row=41; D='6.597'; E='  +0.14%  ' },
    @{ row=42; D='0.8989'; E='  -0.98%  ' },
    @{ row=43; D='0.9998'; E='  -0.20%  ' },
    @{ row=44; D='102.49'; E='  -0.35%  ' },
    @{ row=45; D='1.984.88'; E='  +0.21%  ' },
    @{ row=46; D='66.49'; E='  +2.33%  ' },
    @{ row=47; D='0.00000000122'; E='  -0.79%  ' },
    @{ row=48; D='0.5082'; E='  -0.78%  ' },
    @{ row=49; D='0.4075'; E='  -0.64%  ' },
    @{ row=50; D='9.079'; E='  +0.27%  ' },
    @{ row=51; D='0.05841'; E='  +0.53%  ' }
)

foreach ($item in $updates) {
    $row = $item.row

    if ($item.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $item.D
        $cellD.Style = "Normal"
    }
    if ($item.ContainsKey('E')) {
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $item.E
        $cellE.Style = "Normal"
    }
}
